$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.36733066666667
$ws.Range("H2").Value = 31.101992
$ws.Range("I2").Value = 0.1169328841728879
$ws.Range("J2").Value = 0.1169328841728879
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 136.6884486666667
$ws.Range("N2").Value = 410.065346
$ws.Range("O2").Value = 0.7423691870207686
$ws.Range("P2").Value = 0.7423691870207685
$ws.Range("Q2").Value = 1417.094345641025
$ws.Range("R2").Value = 12753.84911076923
$ws.Range("S2").Value = 0.0868073701594205
$ws.Range("T2").Value = 0.0868073701594205

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.36733066666667
$ws.Range("H3").Value = 31.101992
$ws.Range("I3").Value = 0.1169328841728879
$ws.Range("J3").Value = 0.1169328841728879
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.8952453333333334
$ws.Range("N3").Value = 2.685736
$ws.Range("O3").Value = 0.004862170554817893
$ws.Range("P3").Value = 0.004862170554817893
$ws.Range("Q3").Value = 9.281304398456889
$ws.Range("R3").Value = 83.531739586112
$ws.Range("S3").Value = 0.0005685476263153468
$ws.Range("T3").Value = 0.0005685476263153469

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.36733066666667
$ws.Range("H4").Value = 31.101992
$ws.Range("I4").Value = 0.1169328841728879
$ws.Range("J4").Value = 0.1169328841728879
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 41.63761133333333
$ws.Range("N4").Value = 124.912834
$ws.Range("O4").Value = 0.2261381995079395
$ws.Range("P4").Value = 0.2261381995079395
$ws.Range("Q4").Value = 431.6708848628141
$ws.Range("R4").Value = 3885.037963765328
$ws.Range("S4").Value = 0.0264429918901273
$ws.Range("T4").Value = 0.0264429918901273

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.36733066666667
$ws.Range("H5").Value = 31.101992
$ws.Range("I5").Value = 0.1169328841728879
$ws.Range("J5").Value = 0.1169328841728879
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 4.903320333333333
$ws.Range("N5").Value = 14.709961
$ws.Range("O5").Value = 0.02663044291647413
$ws.Range("P5").Value = 0.02663044291647413
$ws.Range("Q5").Value = 50.83434326025689
$ws.Range("R5").Value = 457.5090893423119
$ws.Range("S5").Value = 0.003113974497024773
$ws.Range("T5").Value = 0.003113974497024773

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 37.91490933333333
$ws.Range("H6").Value = 113.744728
$ws.Range("I6").Value = 0.4276413904453658
$ws.Range("J6").Value = 0.4276413904453659
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 136.6884486666667
$ws.Range("N6").Value = 410.065346
$ws.Range("O6").Value = 0.7423691870207686
$ws.Range("P6").Value = 0.7423691870207685
$ws.Range("Q6").Value = 5182.530138110654
$ws.Range("R6").Value = 46642.77124299589
$ws.Range("S6").Value = 0.3174677913613573
$ws.Range("T6").Value = 0.3174677913613573

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 37.91490933333333
$ws.Range("H7").Value = 113.744728
$ws.Range("I7").Value = 0.4276413904453658
$ws.Range("J7").Value = 0.4276413904453659
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.8952453333333334
$ws.Range("N7").Value = 2.685736
$ws.Range("O7").Value = 0.004862170554817893
$ws.Range("P7").Value = 0.004862170554817893
$ws.Range("Q7").Value = 33.94314564442312
$ws.Range("R7").Value = 305.4883107998081
$ws.Range("S7").Value = 0.002079265376644839
$ws.Range("T7").Value = 0.00207926537664484

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 37.91490933333333
$ws.Range("H8").Value = 113.744728
$ws.Range("I8").Value = 0.4276413904453658
$ws.Range("J8").Value = 0.4276413904453659
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 41.63761133333333
$ws.Range("N8").Value = 124.912834
$ws.Range("O8").Value = 0.2261381995079395
$ws.Range("P8").Value = 0.2261381995079395
$ws.Range("Q8").Value = 1578.686258559906
$ws.Range("R8").Value = 14208.17632703915
$ws.Range("S8").Value = 0.09670605407038675
$ws.Range("T8").Value = 0.09670605407038678

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 37.91490933333333
$ws.Range("H9").Value = 113.744728
$ws.Range("I9").Value = 0.4276413904453658
$ws.Range("J9").Value = 0.4276413904453659
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 4.903320333333333
$ws.Range("N9").Value = 14.709961
$ws.Range("O9").Value = 0.02663044291647413
$ws.Range("P9").Value = 0.02663044291647413
$ws.Range("Q9").Value = 185.9089458706231
$ws.Range("R9").Value = 1673.180512835608
$ws.Range("S9").Value = 0.01138827963697694
$ws.Range("T9").Value = 0.01138827963697694

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 26.72147866666667
$ws.Range("H10").Value = 80.16443599999999
$ws.Range("I10").Value = 0.3013909433702152
$ws.Range("J10").Value = 0.3013909433702153
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 136.6884486666667
$ws.Range("N10").Value = 410.065346
$ws.Range("O10").Value = 0.7423691870207686
$ws.Range("P10").Value = 0.7423691870207685
$ws.Range("Q10").Value = 3652.517465026095
$ws.Range("R10").Value = 32872.65718523485
$ws.Range("S10").Value = 0.2237433496051692
$ws.Range("T10").Value = 0.2237433496051692

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 26.72147866666667
$ws.Range("H11").Value = 80.16443599999999
$ws.Range("I11").Value = 0.3013909433702152
$ws.Range("J11").Value = 0.3013909433702153
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.8952453333333334
$ws.Range("N11").Value = 2.685736
$ws.Range("O11").Value = 0.004862170554817893
$ws.Range("P11").Value = 0.004862170554817893
$ws.Range("Q11").Value = 23.92227907609956
$ws.Range("R11").Value = 215.300511684896
$ws.Range("S11").Value = 0.001465414170343447
$ws.Range("T11").Value = 0.001465414170343448

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 26.72147866666667
$ws.Range("H12").Value = 80.16443599999999
$ws.Range("I12").Value = 0.3013909433702152
$ws.Range("J12").Value = 0.3013909433702153
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 41.63761133333333
$ws.Range("N12").Value = 124.912834
$ws.Range("O12").Value = 0.2261381995079395
$ws.Range("P12").Value = 0.2261381995079395
$ws.Range("Q12").Value = 1112.618542974625
$ws.Range("R12").Value = 10013.56688677162
$ws.Range("S12").Value = 0.0681560052817398
$ws.Range("T12").Value = 0.06815600528173982

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 26.72147866666667
$ws.Range("H13").Value = 80.16443599999999
$ws.Range("I13").Value = 0.3013909433702152
$ws.Range("J13").Value = 0.3013909433702153
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 4.903320333333333
$ws.Range("N13").Value = 14.709961
$ws.Range("O13").Value = 0.02663044291647413
$ws.Range("P13").Value = 0.02663044291647413
$ws.Range("Q13").Value = 131.0239696829996
$ws.Range("R13").Value = 1179.215727146996
$ws.Range("S13").Value = 0.008026174312962804
$ws.Range("T13").Value = 0.008026174312962804

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 13.65680433333333
$ws.Range("H14").Value = 40.970413
$ws.Range("I14").Value = 0.154034782011531
$ws.Range("J14").Value = 0.154034782011531
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 136.6884486666667
$ws.Range("N14").Value = 410.065346
$ws.Range("O14").Value = 0.7423691870207686
$ws.Range("P14").Value = 0.7423691870207685
$ws.Range("Q14").Value = 1866.727398067544
$ws.Range("R14").Value = 16800.5465826079
$ws.Range("S14").Value = 0.1143506758948216
$ws.Range("T14").Value = 0.1143506758948216

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 13.65680433333333
$ws.Range("H15").Value = 40.970413
$ws.Range("I15").Value = 0.154034782011531
$ws.Range("J15").Value = 0.154034782011531
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.8952453333333334
$ws.Range("N15").Value = 2.685736
$ws.Range("O15").Value = 0.004862170554817893
$ws.Range("P15").Value = 0.004862170554817893
$ws.Range("Q15").Value = 12.22619034766311
$ws.Range("R15").Value = 110.035713128968
$ws.Range("S15").Value = 0.0007489433815142589
$ws.Range("T15").Value = 0.0007489433815142591

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 13.65680433333333
$ws.Range("H16").Value = 40.970413
$ws.Range("I16").Value = 0.154034782011531
$ws.Range("J16").Value = 0.154034782011531
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 41.63761133333333
$ws.Range("N16").Value = 124.912834
$ws.Range("O16").Value = 0.2261381995079395
$ws.Range("P16").Value = 0.2261381995079395
$ws.Range("Q16").Value = 568.6367108867157
$ws.Range("R16").Value = 5117.730397980442
$ws.Range("S16").Value = 0.03483314826568557
$ws.Range("T16").Value = 0.03483314826568557

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 13.65680433333333
$ws.Range("H17").Value = 40.970413
$ws.Range("I17").Value = 0.154034782011531
$ws.Range("J17").Value = 0.154034782011531
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 4.903320333333333
$ws.Range("N17").Value = 14.709961
$ws.Range("O17").Value = 0.02663044291647413
$ws.Range("P17").Value = 0.02663044291647413
$ws.Range("Q17").Value = 66.96368637598812
$ws.Range("R17").Value = 602.6731773838929
$ws.Range("S17").Value = 0.004102014469509613
$ws.Range("T17").Value = 0.004102014469509613
